# Add 2022-Q3 data: insert a new quarter sheet ("2022-Q3") before the
# existing "2022-Q2" sheet, populate it with the Q3 fund-holdings data,
# and add a matching summary row on the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: insert a new row for 2022-Q3
#    right after the header, shifting 2022-Q2 / 2022-Q1 down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

# Copy the number-style (centred/bold/bordered) formatting already used
# by column A in the existing data rows onto the new row's A cell.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 0.85

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" worksheet right before "2022-Q2" and
#    fill it with the quarter's fund-holdings table.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Reuse the header/"numbered-row" cell style (bold font, thin border,
# centred) that already exists on the totals sheet so the new sheet's
# look matches its siblings.
$styleSrc = $totalSheet.Range("B1")
$styleSrc.Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2:A8").PasteSpecial(-4122)

$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Columns B-G hold numeric-looking text (fund codes / percentages) that
# must stay as text (leading zeros in fund codes, etc.) - force text
# formatting before writing the values.
$q3Sheet.Range("B2:G8").NumberFormat = "@"

$rows = @(
    @(0, "513300", "华夏纳斯达克100ETF（QDII）",               "11.08", "97.32", "1.74", "0.1928", 5),
    @(1, "000043", "嘉实美国成长股票（QDII）人民币",             "12.41", "92.80", "1.39", "0.1725", 8),
    @(2, "000044", "嘉实美国成长股票（QDII）美元现汇",            "12.41", "92.80", "1.39", "0.1725", 8),
    @(3, "161128", "易方达标普信息科技指数（QDII-LOF）人民币",     "4.99",  "91.96", "2.08", "0.1038", 6),
    @(4, "012868", "易方达标普信息科技指数（QDII-LOF）人民币 C",   "4.99",  "91.96", "2.08", "0.1038", 6),
    @(5, "003721", "易方达标普信息科技指数（QDII-LOF）美元A",      "4.84",  "91.96", "2.08", "0.1007", 6),
    @(6, "012869", "易方达标普信息科技指数（QDII-LOF）美元 C",     "0.15",  "91.96", "2.08", "0.0031", 6)
)

$r = 2
foreach ($row in $rows) {
    $q3Sheet.Range("A$r").Value = $row[0]
    $q3Sheet.Range("B$r").Value = $row[1]
    $q3Sheet.Range("C$r").Value = $row[2]
    $q3Sheet.Range("D$r").Value = $row[3]
    $q3Sheet.Range("E$r").Value = $row[4]
    $q3Sheet.Range("F$r").Value = $row[5]
    $q3Sheet.Range("G$r").Value = $row[6]
    $q3Sheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Drop the temporary text number-format now that the values are typed in
# (keeps them as text cells, just without the leftover "@" style index).
$q3Sheet.Range("B2:G8").ClearFormats()

# ---------------------------------------------------------------------
# 3. Restore the original active sheet ("2022-Q1") since adding the new
#    sheet shifted the selection to it.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
